# Converts a target OOXML EMU coordinate into the "points" value that the
# Shape.Left/Top/Width/Height (float, single precision) COM properties need
# to be set to so that, once the host's points->EMU conversion
# (floor(float32(pts) * 12700)) runs, the stored EMU value matches exactly.
function EmuToPt($emu) {
    return ($emu + 0.5) / 12700
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Refresh the cached "datetimeFigureOut" date field text (10/16/19 ->
#    10/22/19) on the slide master and on every slide layout's Date
#    placeholder.
# ---------------------------------------------------------------------
$newDate = "10/22/19"
$m = $p.SlideMaster

for ($i = 1; $i -le $m.Shapes.Placeholders.Count; $i++) {
    $ph = $m.Shapes.Placeholders.Item($i)
    if ($ph.PlaceholderFormat.Type -eq 16) {
        $ph.TextFrame.TextRange.Text = $newDate
    }
}

for ($li = 1; $li -le $m.CustomLayouts.Count; $li++) {
    $cl = $m.CustomLayouts.Item($li)
    for ($i = 1; $i -le $cl.Shapes.Placeholders.Count; $i++) {
        $ph = $cl.Shapes.Placeholders.Item($i)
        if ($ph.PlaceholderFormat.Type -eq 16) {
            $ph.TextFrame.TextRange.Text = $newDate
        }
    }
}

# ---------------------------------------------------------------------
# 2. On the (single) slide: move the little "e5" label/marker.
#    Duplicate the existing "TextBox 27" (so the new shape keeps its
#    exact run/paragraph formatting), rename + reposition the copy to
#    its new spot, then delete the old oval marker ("Oval 150") and the
#    old "TextBox 27" it was paired with.
# ---------------------------------------------------------------------
$s = $p.Slides.Item(1)

$origLabel = $s.Shapes.Item("TextBox 27")
$newLabel = $origLabel.Duplicate()
$newLabel.Name = "TextBox 3"
$newLabel.Left = EmuToPt(4122624)
$newLabel.Top = EmuToPt(999183)
$newLabel.Width = EmuToPt(338554)
$newLabel.Height = EmuToPt(276999)

$s.Shapes.Item("Oval 150").Delete()
$s.Shapes.Item("TextBox 27").Delete()
